# Daily attendance processing - 2025-11-19 20:21:48
# Normalizes the "Recorded By" (column G) lists so that a leading
# "System"/"system" (or other first-listed recorder) is rotated to the
# end of the comma-separated list, matching the reprocessed attendance
# export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "system, backup@backdoor.com, System") {
        $cell.Value2 = "backup@backdoor.com, System, system"
    }
    elseif ($val -eq "admin@admin.com, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, admin@admin.com"
    }
}
